$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 481, shifting existing rows 481:500 down to 482:501
$ws.Rows(481).Insert()

# Populate the new row 481 with the new record (Ajo / Feria Lagunitas de Puerto Montt)
$ws.Range("A481").Value = 4
$ws.Range("B481").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C481").Value = "Los Lagos"
$ws.Range("D481").Value = 45147
$ws.Range("E481").Value = 10
$ws.Range("F481").Value = 100112003
$ws.Range("G481").Value = "Ajo"
$ws.Range("H481").Value = "Chino"
$ws.Range("I481").Value = "Primera"
$ws.Range("J481").Value = 35
$ws.Range("K481").Value = 23000
$ws.Range("L481").Value = 23000
$ws.Range("M481").Value = 23000
$ws.Range("N481").Value = "$/caja 10 kilos"
$ws.Range("O481").Value = "China"
$ws.Range("P481").Value = 2300
$ws.Range("Q481").Value = 10
$ws.Range("R481").Value = "Hortaliza"
